$wb = $excel.ActiveWorkbook
$wsJournal = $wb.Worksheets.Item("Journal")
$wsTotaux = $wb.Worksheets.Item("Totaux")

# --- Journal sheet: extend the "Tableau1" table from A1:F6 to A1:F7 ---
$table1 = $wsJournal.ListObjects.Item(1)
$table1.Resize($wsJournal.Range("A1:F7")) | Out-Null

# Give the new row 7 the same number formats as the existing data rows (A:D)
# by copying row 6's formatting before filling in values.
$wsJournal.Range("A6:D6").Copy() | Out-Null
$wsJournal.Range("A7:D7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Fill in the existing (previously blank) data rows 2-7 ---

# Row 2: only the "Temps [h]" duration was missing (0:15)
$wsJournal.Range("C2").Value = (15/1440)

# Row 3: 1:00
$wsJournal.Range("A3").Value = 45048
$wsJournal.Range("B3").Value = 1
$wsJournal.Range("C3").Value = (1/24)
$wsJournal.Range("D3").Value = "Documentation"
$wsJournal.Range("E3").Value = "Rédaction de la planification initiale"

# Row 4: 0:45
$wsJournal.Range("A4").Value = 45048
$wsJournal.Range("B4").Value = 1
$wsJournal.Range("C4").Value = (45/1440)
$wsJournal.Range("D4").Value = "Documentation"
$wsJournal.Range("E4").Value = "Rédaction du rapport de projet"
$wsJournal.Range("F4").Value = 'Rédaction de la partie "Analyse préliminaire"'

# Row 5: 1:00
$wsJournal.Range("A5").Value = 45048
$wsJournal.Range("B5").Value = 1
$wsJournal.Range("C5").Value = (1/24)
$wsJournal.Range("D5").Value = "Documentation"
$wsJournal.Range("E5").Value = "Création des différents sprint sur Icescrum"

# Row 6: 1:00
$wsJournal.Range("A6").Value = 45048
$wsJournal.Range("B6").Value = 1
$wsJournal.Range("C6").Value = (1/24)
$wsJournal.Range("D6").Value = "Documentation"
$wsJournal.Range("E6").Value = "Modification de la planification initiale "
$wsJournal.Range("F6").Value = 'J''ai rencontré pas mal de difficulté pour estimé le nombre d''heure de chaque bloc de la planification initiale'

# Row 7 (new): 1:00
$wsJournal.Range("A7").Value = 45048
$wsJournal.Range("B7").Value = 1
$wsJournal.Range("C7").Value = (1/24)
$wsJournal.Range("D7").Value = "Documentation"
$wsJournal.Range("E7").Value = "Rédaction du rapport de projet"
$wsJournal.Range("F7").Value = "Modification de la partie Introduction du rapport"

# --- Totaux sheet: total hours formula, summing the Journal table's column ---
$wsTotaux.Range("B2").Formula = "=SUM(Tableau1[Temps '[h']])"

# --- Restore selections on both sheets (Journal must stay the active tab) ---
$wsTotaux.Range("D12").Select() | Out-Null
$wsJournal.Range("F15").Select() | Out-Null
